# Updates the crypto price/volume table with freshly scraped values.
# Column D ("Price") and column E ("Volume(1h)") hold plain text (not numbers),
# so numeric-looking Price values are forced to Text format before assignment
# -- otherwise Excel auto-coerces strings like "585.09" into a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.337.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.499.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.37%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  +4.93%  "
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.100.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.07%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.505.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.381.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.12%  "
$ws.Range("E19").Value = "  +6.32%  "
$ws.Range("E20").Value = "  +7.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +9.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.642.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("E28").Value = "  +10.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.988"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +6.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +30.87%  "
$ws.Range("E36").Value = "  +5.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "173.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.535.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("E43").Value = "  +8.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +10.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.590.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.74%  "
$ws.Range("E49").Value = "  +12.44%  "
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("E51").Value = "  +5.68%  "
